$wb = $excel.ActiveWorkbook

# Worksheets
$wsStudent = $wb.Worksheets.Item("Event Table (P) Student")
$wsAdmin   = $wb.Worksheets.Item("Event Table (P) Admin")

# --- Content fix: correct the "inforamtion/inforamation" typos on the
# Student sheet (row 7 - "Students updates the profile information" /
# "Updates profile information") ---
$wsStudent.Range("B7").Value = "Students updates the profile information"
$wsStudent.Range("E7").Value = "Updates profile information"

# --- View state: active sheet moves from "Event Table (P) Admin" to
# "Event Table (P) Student", and the selected cell on each sheet changes ---
$wsAdmin.Activate()
$null = $wsAdmin.Range("B5").Select()

$wsStudent.Activate()
$null = $wsStudent.Range("C7").Select()
